$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new day column (17-nov) before the October block ---
$ws = $wb.Worksheets.Item("Prix Spot")

# Inserting the entire column at DV shifts DV:EZ (the October columns) one column
# to the right, to DW:FA, and keeps the header style/border on the new column.
$ws.Range("DV1").EntireColumn.Insert()

# Header for the freshly inserted column.
$ws.Range("DV1").Value = "17-nov"

# The inserted column has no data in rows 2-25 yet; fill with the same
# "no data" placeholder used throughout the rest of the sheet.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 126).Value = "-"
}

# --- Sheet "Gaz": append the two newest daily quotes ---
$wsGaz = $wb.Worksheets.Item("Gaz")

$wsGaz.Range("A153").Formula = "=""2025-11-15"""
$wsGaz.Range("A153").Copy()
$wsGaz.Range("A153").PasteSpecial(-4163)
$wsGaz.Range("B153").Value = 29.305

$wsGaz.Range("A154").Formula = "=""2025-11-16"""
$wsGaz.Range("A154").Copy()
$wsGaz.Range("A154").PasteSpecial(-4163)
$wsGaz.Range("B154").Value = 29.305

# --- Sheet "CO2": append the two newest daily quotes ---
$wsCO2 = $wb.Worksheets.Item("CO2")

$wsCO2.Range("A154").Formula = "=""2025-11-15"""
$wsCO2.Range("A154").Copy()
$wsCO2.Range("A154").PasteSpecial(-4163)
$wsCO2.Range("B154").Value = 80.72

$wsCO2.Range("A155").Formula = "=""2025-11-16"""
$wsCO2.Range("A155").Copy()
$wsCO2.Range("A155").PasteSpecial(-4163)
$wsCO2.Range("B155").Value = 80.72
